$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Trim trailing whitespace from A25 / A26 text values
$ws.Range("A26").Value = "BK 1 CKT INT CI-2/SEL TRBL"
$ws.Range("A25").Value = "LOSS OF POTENTIAL"

# C31/C32 previously held "None" placeholder text -> now "UNDEFINED"
$ws.Range("C31").Value = "UNDEFINED"
$ws.Range("C32").Value = "UNDEFINED"

# F31/F32 previously held "None" placeholder text -> now cleared entirely
$ws.Range("F31").ClearContents()
$ws.Range("F32").ClearContents()

# Move the long note text from column C to column G on rows 34 & 35,
# preserving the wrap-text formatting that was applied to the source cells.
$valC34 = $ws.Range("C34").Value2
$valC35 = $ws.Range("C35").Value2
$ws.Range("C34").Clear()
$ws.Range("C35").Clear()
$ws.Range("G34").Value = $valC34
$ws.Range("G34").WrapText = $true
$ws.Range("G35").Value = $valC35
$ws.Range("G35").WrapText = $true

# Shrink rows 34 & 35 to fit the now-narrower note column
$ws.Rows.Item(34).RowHeight = 135
$ws.Rows.Item(35).RowHeight = 60

# Give column G an explicit width to host the relocated notes
$ws.Columns.Item(7).ColumnWidth = 25

# Update the active selection
$ws.Range("A24").Select()
